$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3275576666666667
$ws.Range("H2").Value = 0.982673
$ws.Range("I2").Value = 0.05486041027915935
$ws.Range("J2").Value = 0.05486041027915935
$ws.Range("M2").Value = 1.123319
$ws.Range("N2").Value = 3.369957
$ws.Range("O2").Value = 0.05053686506648315
$ws.Range("P2").Value = 0.05053686506648315
$ws.Range("Q2").Value = 0.3679517505623334
$ws.Range("R2").Value = 3.311565755061
$ws.Range("S2").Value = 0.002772473151769781
$ws.Range("T2").Value = 0.002772473151769781
$ws.Range("G3").Value = 0.3275576666666667
$ws.Range("H3").Value = 0.982673
$ws.Range("I3").Value = 0.05486041027915935
$ws.Range("J3").Value = 0.05486041027915935
$ws.Range("O3").Value = 0.5042195746532222
$ws.Range("P3").Value = 0.5042195746532223
$ws.Range("Q3").Value = 3.671151246073111
$ws.Range("R3").Value = 33.040361214658
$ws.Range("S3").Value = 0.02766169273625899
$ws.Range("T3").Value = 0.02766169273625899
$ws.Range("G4").Value = 0.3275576666666667
$ws.Range("H4").Value = 0.982673
$ws.Range("I4").Value = 0.05486041027915935
$ws.Range("J4").Value = 0.05486041027915935
$ws.Range("M4").Value = 4.958620666666667
$ws.Range("N4").Value = 14.875862
$ws.Range("O4").Value = 0.2230827962023326
$ws.Range("P4").Value = 0.2230827962023326
$ws.Range("Q4").Value = 1.624234215458445
$ws.Range("R4").Value = 14.618107939126
$ws.Range("S4").Value = 0.01223841372588206
$ws.Range("T4").Value = 0.01223841372588206
$ws.Range("G5").Value = 0.3275576666666667
$ws.Range("H5").Value = 0.982673
$ws.Range("I5").Value = 0.05486041027915935
$ws.Range("J5").Value = 0.05486041027915935
$ws.Range("M5").Value = 4.938126
$ws.Range("N5").Value = 14.814378
$ws.Range("O5").Value = 0.222160764077962
$ws.Range("P5").Value = 0.222160764077962
$ws.Range("Q5").Value = 1.617521030266
$ws.Range("R5").Value = 14.557689272394
$ws.Range("S5").Value = 0.01218783066524852
$ws.Range("T5").Value = 0.01218783066524852
$ws.Range("I6").Value = 0.8684635977749966
$ws.Range("J6").Value = 0.8684635977749967
$ws.Range("M6").Value = 1.123319
$ws.Range("N6").Value = 3.369957
$ws.Range("O6").Value = 0.05053686506648315
$ws.Range("P6").Value = 0.05053686506648315
$ws.Range("Q6").Value = 5.824832506263001
$ws.Range("R6").Value = 52.42349255636701
$ws.Range("S6").Value = 0.0438894276559075
$ws.Range("T6").Value = 0.0438894276559075
$ws.Range("I7").Value = 0.8684635977749966
$ws.Range("J7").Value = 0.8684635977749967
$ws.Range("O7").Value = 0.5042195746532222
$ws.Range("P7").Value = 0.5042195746532223
$ws.Range("S7").Value = 0.4378963458719158
$ws.Range("T7").Value = 0.437896345871916
$ws.Range("I8").Value = 0.8684635977749966
$ws.Range("J8").Value = 0.8684635977749967
$ws.Range("M8").Value = 4.958620666666667
$ws.Range("N8").Value = 14.875862
$ws.Range("O8").Value = 0.2230827962023326
$ws.Range("P8").Value = 0.2230827962023326
$ws.Range("Q8").Value = 25.712317556658
$ws.Range("R8").Value = 231.410858009922
$ws.Range("S8").Value = 0.1937392877915841
$ws.Range("T8").Value = 0.1937392877915841
$ws.Range("I9").Value = 0.8684635977749966
$ws.Range("J9").Value = 0.8684635977749967
$ws.Range("M9").Value = 4.938126
$ws.Range("N9").Value = 14.814378
$ws.Range("O9").Value = 0.222160764077962
$ws.Range("P9").Value = 0.222160764077962
$ws.Range("Q9").Value = 25.606044983502
$ws.Range("R9").Value = 230.454404851518
$ws.Range("S9").Value = 0.1929385364555891
$ws.Range("T9").Value = 0.1929385364555891
$ws.Range("G10").Value = 0.1537743333333333
$ws.Range("H10").Value = 0.461323
$ws.Range("I10").Value = 0.02575461934052592
$ws.Range("J10").Value = 0.02575461934052592
$ws.Range("M10").Value = 1.123319
$ws.Range("N10").Value = 3.369957
$ws.Range("O10").Value = 0.05053686506648315
$ws.Range("P10").Value = 0.05053686506648315
$ws.Range("Q10").Value = 0.1727376303456667
$ws.Range("R10").Value = 1.554638673111
$ws.Range("S10").Value = 0.001301557722450795
$ws.Range("T10").Value = 0.001301557722450795
$ws.Range("G11").Value = 0.1537743333333333
$ws.Range("H11").Value = 0.461323
$ws.Range("I11").Value = 0.02575461934052592
$ws.Range("J11").Value = 0.02575461934052592
$ws.Range("O11").Value = 0.5042195746532222
$ws.Range("P11").Value = 0.5042195746532223
$ws.Range("Q11").Value = 1.723448701950889
$ws.Range("R11").Value = 15.511038317558
$ws.Range("S11").Value = 0.01298598320923563
$ws.Range("T11").Value = 0.01298598320923563
$ws.Range("G12").Value = 0.1537743333333333
$ws.Range("H12").Value = 0.461323
$ws.Range("I12").Value = 0.02575461934052592
$ws.Range("J12").Value = 0.02575461934052592
$ws.Range("M12").Value = 4.958620666666667
$ws.Range("N12").Value = 14.875862
$ws.Range("O12").Value = 0.2230827962023326
$ws.Range("P12").Value = 0.2230827962023326
$ws.Range("Q12").Value = 0.7625085872695557
$ws.Range("R12").Value = 6.862577285426001
$ws.Range("S12").Value = 0.005745412497611197
$ws.Range("T12").Value = 0.005745412497611197
$ws.Range("G13").Value = 0.1537743333333333
$ws.Range("H13").Value = 0.461323
$ws.Range("I13").Value = 0.02575461934052592
$ws.Range("J13").Value = 0.02575461934052592
$ws.Range("M13").Value = 4.938126
$ws.Range("N13").Value = 14.814378
$ws.Range("O13").Value = 0.222160764077962
$ws.Range("P13").Value = 0.222160764077962
$ws.Range("Q13").Value = 0.7593570335660002
$ws.Range("R13").Value = 6.834213302094001
$ws.Range("S13").Value = 0.005721665911228295
$ws.Range("T13").Value = 0.005721665911228296
$ws.Range("G14").Value = 0.3040386666666667
$ws.Range("H14").Value = 0.9121160000000001
$ws.Range("I14").Value = 0.05092137260531806
$ws.Range("J14").Value = 0.05092137260531806
$ws.Range("M14").Value = 1.123319
$ws.Range("N14").Value = 3.369957
$ws.Range("O14").Value = 0.05053686506648315
$ws.Range("P14").Value = 0.05053686506648315
$ws.Range("Q14").Value = 0.3415324110013335
$ws.Range("R14").Value = 3.073791699012001
$ws.Range("S14").Value = 0.00257340653635507
$ws.Range("T14").Value = 0.00257340653635507
$ws.Range("G15").Value = 0.3040386666666667
$ws.Range("H15").Value = 0.9121160000000001
$ws.Range("I15").Value = 0.05092137260531806
$ws.Range("J15").Value = 0.05092137260531806
$ws.Range("O15").Value = 0.5042195746532222
$ws.Range("P15").Value = 0.5042195746532223
$ws.Range("Q15").Value = 3.407558557081778
$ws.Range("R15").Value = 30.668027013736
$ws.Range("S15").Value = 0.02567555283581171
$ws.Range("T15").Value = 0.02567555283581172
$ws.Range("G16").Value = 0.3040386666666667
$ws.Range("H16").Value = 0.9121160000000001
$ws.Range("I16").Value = 0.05092137260531806
$ws.Range("J16").Value = 0.05092137260531806
$ws.Range("M16").Value = 4.958620666666667
$ws.Range("N16").Value = 14.875862
$ws.Range("O16").Value = 0.2230827962023326
$ws.Range("P16").Value = 0.2230827962023326
$ws.Range("Q16").Value = 1.507612415999112
$ws.Range("R16").Value = 13.568511743992
$ws.Range("S16").Value = 0.01135968218725521
$ws.Range("T16").Value = 0.01135968218725521
$ws.Range("G17").Value = 0.3040386666666667
$ws.Range("H17").Value = 0.9121160000000001
$ws.Range("I17").Value = 0.05092137260531806
$ws.Range("J17").Value = 0.05092137260531806
$ws.Range("M17").Value = 4.938126
$ws.Range("N17").Value = 14.814378
$ws.Range("O17").Value = 0.222160764077962
$ws.Range("P17").Value = 0.222160764077962
$ws.Range("Q17").Value = 1.501381244872001
$ws.Range("R17").Value = 13.512431203848
$ws.Range("S17").Value = 0.01131273104589606
$ws.Range("T17").Value = 0.01131273104589606
